$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new combined value for A2: a Python-tuple-style string
# representation of (name, [field1, field2, ...])
$name = "Magister of Worth"
$emDash = [char]0x2014
$fields = @(
    "{4}{W}{B}",
    "Creature $emDash Angel",
    "Flying",
    "Will of the council $emDash When Magister of Worth enters the battlefield, starting with you, each player votes for grace or condemnation. If grace gets more votes, each player returns each creature card from their graveyard to the battlefield. If condemnation gets more votes or the vote is tied, destroy all creatures other than Magister of Worth.",
    "4/4"
)

$fieldList = ($fields | ForEach-Object { "'" + $_ + "'" }) -join ", "
$newValue = "('" + $name + "', [" + $fieldList + "])"

$ws.Range("A2").Value = $newValue

# Remove the now-obsolete rows 3-7 (their content has been folded into A2)
$ws.Rows("3:7").Delete()
